$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.887.00"
$ws.Range("D3").Value = "2.218.26"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'261.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("D6").Value = "'86.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.67%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").Value = "'45.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.84%  "
$ws.Range("D11").Value = "'0.0921"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "'7.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.49%  "
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "2.549.67"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "'14.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "2.215.10"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "43.810.04"
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("D19").Value = "'0.0000103"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "'69.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'2.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.97%  "
$ws.Range("D23").Value = "'231.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "'8.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +5.20%  "
$ws.Range("D27").Value = "'10.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "'39.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.23%  "
$ws.Range("D30").Value = "'2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'174.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").Value = "'20.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "'0.0868"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  +1.90%  "
$ws.Range("E36").Value = "  +3.91%  "
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "'12.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.29%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'63.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.02%  "
$ws.Range("D43").Value = "'5.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.85%  "
$ws.Range("D44").Value = "'0.200"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.19%  "
$ws.Range("D45").Value = "'100.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("E48").Value = "  +4.43%  "
$ws.Range("D49").Value = "'1.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").Value = "'0.436"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("E51").Value = "  +5.60%  "
